$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
# D2: 240X80 PORCELANATO sold by AVILA TORRES RAFAEL ALEJANDRO
$ws1.Range("D2").Value = 472.57
# D5: progress label for that column (now 1 of 3 advisors have sales)
$ws1.Range("D5").Value = "1 de 3"

# --- Sheet 2: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
# F2: junio sales for AVILA TORRES RAFAEL ALEJANDRO
$ws2.Range("F2").Value = 472.57
# F5: junio column total
$ws2.Range("F5").Value = 472.57
# Column F auto-width grew from 11 to 12 characters after the longer value
# (ColumnWidth round-trips with a +5/6 padding offset in this host, so the
# input is pre-compensated to land exactly on the target stored width)
$ws2.Columns.Item(6).ColumnWidth = 12 - 5/6

# --- Sheet 3: "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# Row 2: OTROS group
$ws3.Range("D2").Value = 2183.97
$ws3.Range("E2").Value = -2183.97

# Row 3: PORCELANATO group
$ws3.Range("D3").Value = 48.98
$ws3.Range("E3").Value = 17451.02
$ws3.Range("F3").Value = 0.002798857142857143

# Row 4: TOTAL
$ws3.Range("D4").Value = 2232.95
$ws3.Range("E4").Value = 15267.05
$ws3.Range("F4").Value = 0.1275971428571429

# Column D auto-width grew from 12 to 13 characters after the longer value
$ws3.Columns.Item(4).ColumnWidth = 13 - 5/6
